$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Move the "PWM" label from K15 (next to GPIO26) down to K22 (next to GPIO12),
# matching the formatting already used by the other K-column legend entries.
$pwmValue = $ws.Range("K15").Value()
$ws.Range("K15").Clear()
$ws.Range("K22").Value() = $pwmValue
$ws.Range("K18").Copy()
$ws.Range("K22").PasteSpecial($xlPasteFormats)

# Add a new "wire color" column (L) alongside the existing K-column legend,
# copying the formatting from the neighboring K cell in each row.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial($xlPasteFormats)
$ws.Range("L4").Value() = "white"

$ws.Range("K17").Copy()
$ws.Range("L17").PasteSpecial($xlPasteFormats)
$ws.Range("L17").Value() = "green"

$ws.Range("K18").Copy()
$ws.Range("L18").PasteSpecial($xlPasteFormats)
$ws.Range("L18").Value() = "green"

$ws.Range("K22").Copy()
$ws.Range("L22").PasteSpecial($xlPasteFormats)
$ws.Range("L22").Value() = "purple"

$ws.Range("K24").Copy()
$ws.Range("L24").PasteSpecial($xlPasteFormats)
$ws.Range("L24").Value() = "yellow"

$ws.Range("K25").Copy()
$ws.Range("L25").PasteSpecial($xlPasteFormats)
$ws.Range("L25").Value() = "yellow"

# Extra styled (empty) cell next to the moved PWM row, matching the new column band
$ws.Range("K22").Copy()
$ws.Range("M22").PasteSpecial($xlPasteFormats)
$ws.Range("M22").ClearContents()

$excel.CutCopyMode = $false

# Extend column K's width to also cover the new L column, and size the new M column
$ws.Columns.Item(11).ColumnWidth = 8.7265625
$ws.Columns.Item(12).ColumnWidth = 8.7265625
$ws.Columns.Item(13).ColumnWidth = 9.54296875

# Reflect the scrolled/selected view from the edit session
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("L6").Select()
